# Der erste Punkt des Temperaturprofils ist immer bei Raumtemperatur
#
# - Remove the first data point (0.0 min / 20.0 degC, "room temperature"),
#   shifting the remaining profile points up by one row.
# - Add a new "Beschreibung" column describing every remaining phase of
#   the temperature profile.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The profile always starts at room temperature now, so drop that row
# (everything below moves up one row automatically).
$ws.Range("A2:B2").Delete(-4162)

# New descriptive column next to the timestamp/temperature columns.
$ws.Range("C1").Value = "Beschreibung"
$ws.Range("C2").Value = "Aufheizen 1"
$ws.Range("C3").Value = "Aufheizen 2"
$ws.Range("C4").Value = "Aufheizen 3"
$ws.Range("C5").Value = "Halten"
$ws.Range("C6").Value = "Abkühlen 1"
$ws.Range("C7").Value = "Abkühlen 2"

# Match the formatting already used by the other two columns.
$ws.Range("B1:B7").Copy()
$ws.Range("C1:C7").PasteSpecial(-4122)
